# Apply the "Excel Formulas and Functions" practice-workbook edit:
#  - define the `taxrate` name (-> 'Tax Rate'!$A$2)
#  - fill in the Totals / Tax / Profit Split formulas on "Formula Fundamentals"
#  - fill in the missing Tax formula on "Order of Operations"
#  - fill in the Sum / 2nd-most-sold formulas on "Basic Functions"
#  - restore each sheet's last-used selection/zoom and make
#    "Conditional Functions" the active sheet/tab

$wb = $excel.ActiveWorkbook

# ---- Workbook-level defined name --------------------------------------
# (Note: PowerShell double-quoted strings interpolate `$`, so the two
#  dollar signs of the absolute reference must be escaped with backticks.)
$wb.Names.Add("taxrate", "='Tax Rate'!`$A`$2")

# ---- "Formula Fundamentals" --------------------------------------------
$wsFormula = $wb.Worksheets.Item("Formula Fundamentals")

$wsFormula.Range("E2").Formula = "=B2+C2+D2"
$wsFormula.Range("E3").Formula = "=B3+C3+D3"

$wsFormula.Range("B4").Formula = "=B2-B3"
$wsFormula.Range("C4").Formula = "=C2-C3"
$wsFormula.Range("D4").Formula = "=D2-D3"
$wsFormula.Range("E4").Formula = "=E2-E3"

$wsFormula.Range("B5").Formula = "=B4*taxrate"
$wsFormula.Range("C5").Formula = "=C4*taxrate"
$wsFormula.Range("D5").Formula = "=D4*taxrate"
$wsFormula.Range("E5").Formula = "=E4*taxrate"

$wsFormula.Range("B9").Formula = "=E4/2"

# ---- "Order of Operations" ---------------------------------------------
$wsOrder = $wb.Worksheets.Item("Order of Operations")
$wsOrder.Range("B5").Formula = "=B4*D2"

# ---- "Basic Functions" --------------------------------------------------
$wsBasic = $wb.Worksheets.Item("Basic Functions")
$wsBasic.Range("E2").Formula = "=SUM(B2:B8)"
$wsBasic.Range("E5").Formula = "=LARGE(B2:B8,2)"

# ---- Restore per-sheet selection / zoom, matching the saved workbook ----
$wsFormula.Select()
$wsFormula.Range("B10").Select()
$excel.ActiveWindow.Zoom = 122

$wsTaxRate = $wb.Worksheets.Item("Tax Rate")
$wsTaxRate.Select()
$wsTaxRate.Range("A2").Select()

$wsOrder.Select()
$wsOrder.Range("B6").Select()

$wsBasic.Select()
$wsBasic.Range("E6").Select()

$wsBonus = $wb.Worksheets.Item("Bonus")
$wsBonus.Select()
$excel.ActiveWindow.Zoom = 239

# "Conditional Functions" is the sheet that ends up active/selected.
$wsCond = $wb.Worksheets.Item("Conditional Functions")
$wsCond.Select()
